$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.518.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.72%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.587.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.98%  "

$ws.Range("E4").Value = "  +0.99%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "

$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("E7").Value = "  +1.00%  "

$ws.Range("E8").Value = "  +5.51%  "

$ws.Range("E9").Value = "  +0.55%  "

$ws.Range("E10").Value = "  +0.97%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0886"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.813.88"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.94%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.585.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.28%  "

$ws.Range("E14").Value = "  +1.77%  "

$ws.Range("E15").Value = "  -0.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.527.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.75%  "

$ws.Range("E17").Value = "  +1.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.13%  "

$ws.Range("E20").Value = "  +0.34%  "

$ws.Range("E21").Value = "  +0.98%  "

$ws.Range("E22").Value = "  -1.69%  "

$ws.Range("E23").Value = "  -0.93%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.72%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.02%  "

$ws.Range("E26").Value = "  +0.41%  "

$ws.Range("E27").Value = "  -0.92%  "

$ws.Range("E28").Value = "  -0.75%  "

$ws.Range("E29").Value = "  +0.96%  "

$ws.Range("E30").Value = "  -0.85%  "

$ws.Range("E31").Value = "  -0.54%  "

$ws.Range("E32").Value = "  +0.25%  "

$ws.Range("E33").Value = "  +1.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.391.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.48%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.38%  "

$ws.Range("E36").Value = "  -10.47%  "

$ws.Range("E37").Value = "  +1.13%  "

$ws.Range("E38").Value = "  +10.90%  "

$ws.Range("E39").Value = "  -0.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.540"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.811"
$ws.Range("D41").Style = "Normal"

$ws.Range("E42").Value = "  +1.00%  "

$ws.Range("E43").Value = "  -0.42%  "

$ws.Range("E44").Value = "  +0.51%  "

$ws.Range("E45").Value = "  +0.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "62.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.724.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.90%  "

$ws.Range("E48").Value = "  +1.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("E50").Value = "  +0.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0521"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.13%  "
